# Applies the social_media_graph5_data.xlsx edit:
#  - Renames the "indic_is" (column B) category codes from the old verb
#    names to the new E_SM_* codes used by Eurostat.
#  - Leaves columns A, C, D, E, F untouched.
#  - Moves the active selection to L8.
#  - Widens column B to fit the longer codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old label -> new code mapping (column B, rows 2-25)
$map = @{
    "obtain"      = "E_SM_PCUQOR"
    "develop"     = "E_SM_PADVERT"
    "recruit"     = "E_SM_PRCR"
    "exchange"    = "E_SM_PEXCHVOC"
    "involve"     = "E_SM_PCUDEV"
    "collaborate" = "E_SM_PBPCOLL"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Value2
    if ($map.ContainsKey($old)) {
        $cell.Value = $map[$old]
    }
}

# Widen column B (bestFit) to fit the new, longer codes - mirrors Excel
# silently re-fitting an already best-fit column when its contents change.
$ws.Columns.Item(2).AutoFit() | Out-Null

# Move the active selection to L8.
$ws.Range("L8").Select()
